# This workbook contains duplicate-looking rows (same product name / price
# point appearing on consecutive lines). The source data feed had the rows
# in the wrong order relative to their batch/expiry numbers (column B) and
# the dependent Sales/Value figures (columns D-G). This swaps each group of
# mis-ordered rows back into the correct order while leaving the serial
# number (column A) and the trailing blank columns (H-M) untouched.
#
# For a group of rows, row[i] receives the B..G values that currently sit
# in row[i+1] (wrapping around for the last row in the group) -- i.e. a
# cyclic rotation. For 2-row groups this is just a swap.

function Rotate-RowValues {
    param(
        $ws,
        $cols,
        $rows
    )

    $n = $rows.Count

    foreach ($col in $cols) {
        # Snapshot the current values for this column across all rows in the group.
        $orig = @()
        for ($i = 0; $i -lt $n; $i++) {
            $orig += $ws.Range($col + $rows[$i]).Value2
        }

        # Write back the rotated values: row[i] <- row[i+1] (wrap-around).
        for ($i = 0; $i -lt $n; $i++) {
            $nextIndex = ($i + 1) % $n
            $ws.Range($col + $rows[$i]).Value2 = $orig[$nextIndex]
        }
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-row data which needs reordering.
$cols = @("B", "C", "D", "E", "F", "G")

# Each entry is a group of worksheet rows (by row number) whose B..G data
# needs to be cyclically rotated into the correct order.
$rowGroups = @(
    , @(264, 265)
    , @(346, 347)
    , @(350, 351, 352)
    , @(375, 376)
    , @(382, 383)
    , @(389, 390)
    , @(419, 420)
    , @(431, 432)
    , @(536, 537)
    , @(583, 584)
    , @(586, 587)
    , @(590, 591)
    , @(601, 602)
    , @(687, 688)
    , @(709, 710)
    , @(720, 721)
    , @(859, 860)
    , @(889, 890)
)

foreach ($group in $rowGroups) {
    Rotate-RowValues $ws $cols $group
}
